# Refined metadata to be additional tab
#
# 1) Update the "panel_query_time" (column F) timestamps on the "data" sheet
#    to reflect the later query run.
# 2) Add a new "metadata" worksheet (after "data") summarising the panel
#    query that produced the data sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1) refresh the per-row query timestamps on "data" -----------------
$dataSheet.Range("F2").Value  = "2021-10-05 14:22:07.845430"
$dataSheet.Range("F3").Value  = "2021-10-05 14:22:07.845438"
$dataSheet.Range("F4").Value  = "2021-10-05 14:22:07.845441"
$dataSheet.Range("F5").Value  = "2021-10-05 14:22:07.845444"
$dataSheet.Range("F6").Value  = "2021-10-05 14:22:07.845447"
$dataSheet.Range("F7").Value  = "2021-10-05 14:22:07.845449"
$dataSheet.Range("F8").Value  = "2021-10-05 14:22:07.845452"
$dataSheet.Range("F9").Value  = "2021-10-05 14:22:07.845454"
$dataSheet.Range("F10").Value = "2021-10-05 14:22:07.845457"
$dataSheet.Range("F11").Value = "2021-10-05 14:22:07.845459"
$dataSheet.Range("F12").Value = "2021-10-05 14:22:07.845462"
$dataSheet.Range("F13").Value = "2021-10-05 14:22:07.845464"
$dataSheet.Range("F14").Value = "2021-10-05 14:22:07.845467"
$dataSheet.Range("F15").Value = "2021-10-05 14:22:07.845469"
$dataSheet.Range("F16").Value = "2021-10-05 14:22:07.845472"
$dataSheet.Range("F17").Value = "2021-10-05 14:22:07.845474"
$dataSheet.Range("F18").Value = "2021-10-05 14:22:07.845477"
$dataSheet.Range("F19").Value = "2021-10-05 14:22:07.845480"
$dataSheet.Range("F20").Value = "2021-10-05 14:22:07.845482"
$dataSheet.Range("F21").Value = "2021-10-05 14:22:07.845485"
$dataSheet.Range("F22").Value = "2021-10-05 14:22:07.845487"
$dataSheet.Range("F23").Value = "2021-10-05 14:22:07.845489"
$dataSheet.Range("F24").Value = "2021-10-05 14:22:07.845492"
$dataSheet.Range("F25").Value = "2021-10-05 14:22:07.845494"
$dataSheet.Range("F26").Value = "2021-10-05 14:22:07.845497"
$dataSheet.Range("F27").Value = "2021-10-05 14:22:07.845500"
$dataSheet.Range("F28").Value = "2021-10-05 14:22:07.845502"
$dataSheet.Range("F29").Value = "2021-10-05 14:22:07.845504"
$dataSheet.Range("F30").Value = "2021-10-05 14:22:07.845507"
$dataSheet.Range("F31").Value = "2021-10-05 14:22:07.845509"
$dataSheet.Range("F32").Value = "2021-10-05 14:22:07.845512"
$dataSheet.Range("F33").Value = "2021-10-05 14:22:07.845514"
$dataSheet.Range("F34").Value = "2021-10-05 14:22:07.845517"
$dataSheet.Range("F35").Value = "2021-10-05 14:22:07.845520"
$dataSheet.Range("F36").Value = "2021-10-05 14:22:07.845522"
$dataSheet.Range("F37").Value = "2021-10-05 14:22:07.845525"
$dataSheet.Range("F38").Value = "2021-10-05 14:22:07.845527"
$dataSheet.Range("F39").Value = "2021-10-05 14:22:07.845530"
$dataSheet.Range("F40").Value = "2021-10-05 14:22:07.845532"
$dataSheet.Range("F41").Value = "2021-10-05 14:22:07.845534"
$dataSheet.Range("F42").Value = "2021-10-05 14:22:07.845537"
$dataSheet.Range("F43").Value = "2021-10-05 14:22:07.845540"
$dataSheet.Range("F44").Value = "2021-10-05 14:22:07.845542"
$dataSheet.Range("F45").Value = "2021-10-05 14:22:07.845545"
$dataSheet.Range("F46").Value = "2021-10-05 14:22:07.845547"
$dataSheet.Range("F47").Value = "2021-10-05 14:22:07.845550"

# --- 2) add the "metadata" sheet right after "data" ---------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Palmoplantar keratoderma and erythrokeratodermas"
$ws.Range("C2").Value = 215

# data_version ("1.20") must stay textual (not collapse to the number 1.2)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.20"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2021-08-31T15:56:19.349214Z"
$ws.Range("F2").Value = "2021-10-05 14:22:07.841813"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/215/?format=json"

# match the "data" sheet's header/index styling (bold, centered, bordered)
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# keep "data" as the active/selected sheet
$dataSheet.Activate()
